$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.29248327091868
$ws.Range("C2").Value = 12.58150029580305
$ws.Range("D2").Value = 14.94586618187336
$ws.Range("E2").Value = 16.35877698076819
$ws.Range("G2").Value = 3.657733676325439
$ws.Range("I2").Value = 21.98034412123139
$ws.Range("J2").Value = 9.335429618785826
$ws.Range("O2").Value = 25.71922603586754
$ws.Range("B3").Value = 16.62660755172334
$ws.Range("C3").Value = 12.00174920082771
$ws.Range("D3").Value = 14.88487520926003
$ws.Range("E3").Value = 16.29818104252226
$ws.Range("G3").Value = 3.660509764940644
$ws.Range("I3").Value = 22.14816042057675
$ws.Range("J3").Value = 9.344455397663488
$ws.Range("O3").Value = 25.82305430162294
$ws.Range("B4").Value = 16.20505375414281
$ws.Range("C4").Value = 11.6317947390259
$ws.Range("D4").Value = 14.85084665927451
$ws.Range("E4").Value = 16.26475057894311
$ws.Range("G4").Value = 3.662302824207698
$ws.Range("I4").Value = 22.25758186391933
$ws.Range("J4").Value = 9.351474296167506
$ws.Range("O4").Value = 25.89504089683061
$ws.Range("B5").Value = 16.03031354382879
$ws.Range("C5").Value = 11.4777051916604
$ws.Range("D5").Value = 14.83784945556937
$ws.Range("E5").Value = 16.25208590574393
$ws.Range("G5").Value = 3.66305584955458
$ws.Range("I5").Value = 22.30377384070133
$ws.Range("J5").Value = 9.354705900255068
$ws.Range("O5").Value = 25.92643493325519
$ws.Range("B6").Value = 16.00112724436167
$ws.Range("C6").Value = 11.45192361437634
$ws.Range("D6").Value = 14.83574406851017
$ws.Range("E6").Value = 16.25004108027465
$ws.Range("G6").Value = 3.663182240349971
$ws.Range("I6").Value = 22.31154065136106
$ws.Range("J6").Value = 9.355264930558228
$ws.Range("O6").Value = 25.93177188520325
$ws.Range("B7").Value = 16.20270876755367
$ws.Range("C7").Value = 11.62972984986021
$ws.Range("D7").Value = 14.85066784108841
$ws.Range("E7").Value = 16.26457588663441
$ws.Range("G7").Value = 3.662312889218045
$ws.Range("I7").Value = 22.25819834295054
$ws.Range("J7").Value = 9.351516375294519
$ws.Range("O7").Value = 25.89545596751487
$ws.Range("B8").Value = 17.06566204763628
$ws.Range("C8").Value = 12.384608155216
$ws.Range("D8").Value = 14.92413278850119
$ws.Range("E8").Value = 16.33710566192127
$ws.Range("G8").Value = 3.658672544189854
$ws.Range("I8").Value = 22.03688093776935
$ws.Range("J8").Value = 9.338235130193183
$ws.Range("O8").Value = 25.75330960409177
$ws.Range("B9").Value = 18.64787639085223
$ws.Range("C9").Value = 13.747094482816
$ws.Range("D9").Value = 15.09485854075092
$ws.Range("E9").Value = 16.50884262694824
$ws.Range("G9").Value = 3.652232734856398
$ws.Range("I9").Value = 21.65364489963336
$ws.Range("J9").Value = 9.323911158013463
$ws.Range("O9").Value = 25.54042825810741
$ws.Range("B10").Value = 19.73261872113134
$ws.Range("C10").Value = 14.66871907678343
$ws.Range("D10").Value = 15.23583597741123
$ws.Range("E10").Value = 16.65232616735073
$ws.Range("G10").Value = 3.6479225021898
$ws.Range("I10").Value = 21.40321890588839
$ws.Range("J10").Value = 9.320530652873002
$ws.Range("O10").Value = 25.42486783426689
$ws.Range("B11").Value = 20.20742692627833
$ws.Range("C11").Value = 15.06955078148175
$ws.Range("D11").Value = 15.303174031553
$ws.Range("E11").Value = 16.72119179026204
$ws.Range("G11").Value = 3.646052042221716
$ws.Range("I11").Value = 21.29610008805856
$ws.Range("J11").Value = 9.320541819913313
$ws.Range("O11").Value = 25.38130392910395
$ws.Range("B12").Value = 20.38441982884556
$ws.Range("C12").Value = 15.21860808323015
$ws.Range("D12").Value = 15.32911841946641
$ws.Range("E12").Value = 16.74777060096919
$ws.Range("G12").Value = 3.645356649912709
$ws.Range("I12").Value = 21.25651936309073
$ws.Range("J12").Value = 9.320768430829936
$ws.Range("O12").Value = 25.36611297818122
$ws.Range("B13").Value = 20.34642787758637
$ws.Range("C13").Value = 15.18662842928906
$ws.Range("D13").Value = 15.3235113094289
$ws.Range("E13").Value = 16.74202436465605
$ws.Range("G13").Value = 3.645505842054813
$ws.Range("I13").Value = 21.26499998807578
$ws.Range("J13").Value = 9.320709742324722
$ws.Range("O13").Value = 25.36932638085261
$ws.Range("B14").Value = 20.22204514420692
$ws.Range("C14").Value = 15.08186891844127
$ws.Range("D14").Value = 15.30529967056125
$ws.Range("E14").Value = 16.72336850091619
$ws.Range("G14").Value = 3.645994573552299
$ws.Range("I14").Value = 21.29282403160899
$ws.Range("J14").Value = 9.320556008850014
$ws.Range("O14").Value = 25.38002793108107
$ws.Range("B15").Value = 20.14548819894413
$ws.Range("C15").Value = 15.01734315707313
$ws.Range("D15").Value = 15.29420196424681
$ws.Range("E15").Value = 16.71200598805225
$ws.Range("G15").Value = 3.646295614852886
$ws.Range("I15").Value = 21.30999520553107
$ws.Range("J15").Value = 9.320490791125916
$ws.Range("O15").Value = 25.38675329947175
$ws.Range("B16").Value = 19.70120255456836
$ws.Range("C16").Value = 14.6421459488416
$ws.Range("D16").Value = 15.23149843905149
$ws.Range("E16").Value = 16.64789661928995
$ws.Range("G16").Value = 3.648046551709779
$ws.Range("I16").Value = 21.41035663805035
$ws.Range("J16").Value = 9.320561063935182
$ws.Range("O16").Value = 25.42789698195424
$ws.Range("B17").Value = 19.42377701338887
$ws.Range("C17").Value = 14.40719522770491
$ws.Range("D17").Value = 15.19384191665813
$ws.Range("E17").Value = 16.60947711049441
$ws.Range("G17").Value = 3.649143767365882
$ws.Range("I17").Value = 21.47367044073223
$ws.Range("J17").Value = 9.32100067870169
$ws.Range("O17").Value = 25.45545160005907
$ws.Range("B18").Value = 19.26246037351056
$ws.Range("C18").Value = 14.27032784026053
$ws.Range("D18").Value = 15.17248573603653
$ws.Range("E18").Value = 16.58771854743078
$ws.Range("G18").Value = 3.649783359305808
$ws.Range("I18").Value = 21.51072672251081
$ws.Range("J18").Value = 9.321399357837509
$ws.Range("O18").Value = 25.47214727632486
$ws.Range("B19").Value = 19.20754508307496
$ws.Range("C19").Value = 14.22369245036669
$ws.Range("D19").Value = 15.16530739566928
$ws.Range("E19").Value = 16.58041022184942
$ws.Range("G19").Value = 3.65000137668833
$ws.Range("I19").Value = 21.5233830959201
$ws.Range("J19").Value = 9.321559398599671
$ws.Range("O19").Value = 25.47794527784745
$ws.Range("B20").Value = 19.45349139452215
$ws.Range("C20").Value = 14.43238580130261
$ws.Range("D20").Value = 15.19781928808166
$ws.Range("E20").Value = 16.61353192690682
$ws.Range("G20").Value = 3.649026087436216
$ws.Range("I20").Value = 21.46686432126664
$ws.Range("J20").Value = 9.320938790981764
$ws.Range("O20").Value = 25.45243062350667
$ws.Range("B21").Value = 20.25865645206681
$ws.Range("C21").Value = 15.11271396050524
$ws.Range("D21").Value = 15.31063692818762
$ws.Range("E21").Value = 16.72883471778551
$ws.Range("G21").Value = 3.645850671461165
$ws.Range("I21").Value = 21.28462472147218
$ws.Range("J21").Value = 9.32059513203666
$ws.Range("O21").Value = 25.37684910496658
$ws.Range("B22").Value = 20.76847282771688
$ws.Range("C22").Value = 15.54141335503414
$ws.Range("D22").Value = 15.38695509614587
$ws.Range("E22").Value = 16.80710333544812
$ws.Range("G22").Value = 3.643850571094881
$ws.Range("I22").Value = 21.17125182951938
$ws.Range("J22").Value = 9.321666511332786
$ws.Range("O22").Value = 25.33506777982641
$ws.Range("B23").Value = 20.49791230609589
$ws.Range("C23").Value = 15.31408944316744
$ws.Range("D23").Value = 15.34599165122158
$ws.Range("E23").Value = 16.76506898939043
$ws.Range("G23").Value = 3.644911203841045
$ws.Range("I23").Value = 21.23123499917247
$ws.Range("J23").Value = 9.320976256738726
$ws.Range("O23").Value = 25.35666699988249
$ws.Range("B24").Value = 19.44006318512315
$ws.Range("C24").Value = 14.42100271568232
$ws.Range("D24").Value = 15.19602020441381
$ws.Range("E24").Value = 16.61169771741536
$ws.Range("G24").Value = 3.649079263151175
$ws.Range("I24").Value = 21.46993932365767
$ws.Range("J24").Value = 9.320966315804375
$ws.Range("O24").Value = 25.45379374704147
$ws.Range("B25").Value = 18.23281359871381
$ws.Range("C25").Value = 13.39198239576831
$ws.Range("D25").Value = 15.04588743553761
$ws.Range("E25").Value = 16.45929013211867
$ws.Range("G25").Value = 3.653900567625417
$ws.Range("I25").Value = 21.75186686599839
$ws.Range("J25").Value = 9.326531251081432
$ws.Range("O25").Value = 25.59089396394387
